$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# "Name" row's value is cleared, and its former value is moved up as the
# new "Title" row's value (the old "Title" value is discarded entirely).
$ws1.Range("B4").Value = ""
$ws1.Range("B5").Value = "Mapping Métier/CDA/FHIR : `"Types des résultats classés par type d’examens (BIO, IMG, etc…)`""

# Date value bump
$ws1.Range("B8").Value = "2026-01-07T15:20:53+00:00"

# --- Mapping Table 0 sheet ---
$ws2 = $wb.Worksheets.Item("Mapping Table 0")
$ws2.Range("D10").Value = "FRCDAResultats.component:frResultat"

# --- Mapping Table 1 sheet ---
$ws3 = $wb.Worksheets.Item("Mapping Table 1")
$ws3.Range("D7").Value = "FRDiagnosticReportDocument.performer.extension:performerFunction"
$ws3.Range("D8").Value = "FRDiagnosticReportDocument.resultsInterpreter.extension:performerFunction"
$ws3.Range("A9").Value = "FRCDAResultats.component:frResultat"
